# Added "Locale:Armenian" column to the existing "Locale: English(Uganda)" data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data.
$ws.Rows.Item(1).Insert()

# B1 is written first so its text becomes shared-string #0 (Armenian column),
# then A1 becomes shared-string #1 (the pre-existing Uganda locale column).
$ws.Range("B1").Value = "Locale:Armenian"
$ws.Range("A1").Value = "Locale: English(Uganda)"

# Populate the new "Locale:Armenian" date column (column B).
$armenianDates = @{
    2  = "2/8/2016"
    3  = "12/5/2016"
    4  = "10/19/2015"
    5  = "10/12/2017"
    6  = "1/30/2021"
    7  = "8/1/2010"
    8  = "6/11/2018"
    9  = "6/25/2018"
    10 = "6/23/2017"
    11 = "5/31/2018"
    12 = "4/10/2019"
    13 = "1/30/2019"
    14 = "8/20/2018"
    15 = "2/15/2016"
    16 = "1/1/2014"
    17 = "6/29/2018"
    18 = "8/20/2018"
    19 = "4/28/2017"
    20 = "4/28/2017"
    21 = "2/13/2019"
    22 = "6/26/2018"
    23 = "3/7/2018"
    24 = "7/23/2018"
    25 = "2/9/2018"
    26 = "9/24/2018"
    27 = "6/17/2019"
    28 = "9/24/2018"
    29 = "6/17/2019"
}
foreach ($row in $armenianDates.Keys) {
    $ws.Cells.Item($row, 2).Value = $armenianDates[$row]
}

# Format the new column: Comma base style, custom date format, small red font and a thin box border.
$dataRange = $ws.Range("B2:B29")
$dataRange.Style = "Comma"
$dataRange.NumberFormat = "dd/mm/yyyy;@"
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 10
$dataRange.Font.Color = 255
$dataRange.Borders.LineStyle = 1

# A couple of rows were highlighted yellow by the author.
$ws.Range("B5:B6").Interior.Color = 65535

# One row (B22) was formatted with plain right alignment instead of the Comma look.
$ws.Range("B22").HorizontalAlignment = -4152

# Widen the new column to fit the locale text.
$ws.Columns.Item(2).ColumnWidth = 35

# Re-anchor the existing conditional formatting rule onto the shifted data range.
$cf = $ws.Range("A1:A28").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A2:A29"))

# Restore the active selection as left by the author.
$ws.Range("E5").Select()
